$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2033426183844011
$ws.Range("C2").Value = 0.5153203342618384
$ws.Range("J2").Value = 0.008356545961002786
$ws.Range("P2").Value = 0.1782729805013928
$ws.Range("S2").Value = 0.0947075208913649

# Row 3
$ws.Range("C3").Value = 0.005128205128205128
$ws.Range("J3").Value = 0.05128205128205128
$ws.Range("P3").Value = 0.7589743589743589
$ws.Range("S3").Value = 0.1846153846153846

# Row 4
$ws.Range("J4").Value = 0.119047619047619
$ws.Range("P4").Value = 0.6190476190476191
$ws.Range("S4").Value = 0.2619047619047619

# Row 6
$ws.Range("B6").Value = 0.01834862385321101
$ws.Range("D6").Value = 0.01376146788990826
$ws.Range("F6").Value = 0.06422018348623854
$ws.Range("J6").Value = 0.3256880733944954
$ws.Range("O6").Value = 0.03211009174311927
$ws.Range("Q6").Value = 0.1559633027522936
$ws.Range("R6").Value = 0.05045871559633028
$ws.Range("S6").Value = 0.3394495412844037

# Row 7
$ws.Range("B7").Value = 0.1546961325966851
$ws.Range("D7").Value = 0.005524861878453038
$ws.Range("F7").Value = 0.04972375690607735
$ws.Range("J7").Value = 0.1325966850828729
$ws.Range("O7").Value = 0.05524861878453038
$ws.Range("Q7").Value = 0.1546961325966851
$ws.Range("R7").Value = 0.04972375690607735
$ws.Range("S7").Value = 0.3977900552486188

# Row 8
$ws.Range("B8").Value = 0.09420289855072464
$ws.Range("D8").Value = 0.01932367149758454
$ws.Range("E8").Value = 0.002415458937198068
$ws.Range("F8").Value = 0.05072463768115942
$ws.Range("J8").Value = 0.1231884057971015
$ws.Range("O8").Value = 0.03381642512077294
$ws.Range("Q8").Value = 0.1714975845410628
$ws.Range("R8").Value = 0.1014492753623188
$ws.Range("S8").Value = 0.4033816425120773

# Row 9
$ws.Range("B9").Value = 0.1128205128205128
$ws.Range("D9").Value = 0.02564102564102564
$ws.Range("E9").Value = 0.005128205128205128
$ws.Range("F9").Value = 0.04615384615384616
$ws.Range("J9").Value = 0.1333333333333333
$ws.Range("O9").Value = 0.01025641025641026
$ws.Range("Q9").Value = 0.1230769230769231
$ws.Range("R9").Value = 0.1076923076923077
$ws.Range("S9").Value = 0.4358974358974359

# Row 10
$ws.Range("B10").Value = 0.1276595744680851
$ws.Range("D10").Value = 0.01907556859867938
$ws.Range("E10").Value = 0.0007336757153338225
$ws.Range("F10").Value = 0.06162876008804109
$ws.Range("J10").Value = 0.136463683052091
$ws.Range("O10").Value = 0.02494497432134996
$ws.Range("Q10").Value = 0.2010271460014673
$ws.Range("R10").Value = 0.07703595011005136
$ws.Range("S10").Value = 0.3514306676449009

# Row 11
$ws.Range("G11").Value = 0.1888111888111888
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.2307692307692308
$ws.Range("L11").Value = 0.493006993006993
$ws.Range("S11").Value = 0.01048951048951049

# Row 12
$ws.Range("G12").Value = 0.7310344827586207
$ws.Range("J12").Value = 0.2206896551724138
$ws.Range("K12").Value = 0.01379310344827586
$ws.Range("L12").Value = 0.01379310344827586
$ws.Range("S12").Value = 0.02068965517241379

# Row 13
$ws.Range("G13").Value = 0.6590909090909091
$ws.Range("J13").Value = 0.2272727272727273
$ws.Range("S13").Value = 0.1136363636363636

# Row 15
$ws.Range("F15").Value = 0.04149377593360996
$ws.Range("H15").Value = 0.1535269709543569
$ws.Range("I15").Value = 0.05809128630705394
$ws.Range("J15").Value = 0.3568464730290457
$ws.Range("K15").Value = 0.04149377593360996
$ws.Range("M15").Value = 0.008298755186721992
$ws.Range("O15").Value = 0.05809128630705394
$ws.Range("S15").Value = 0.2821576763485477

# Row 16
$ws.Range("F16").Value = 0.02575107296137339
$ws.Range("H16").Value = 0.1373390557939914
$ws.Range("I16").Value = 0.07296137339055794
$ws.Range("J16").Value = 0.3905579399141631
$ws.Range("K16").Value = 0.1158798283261803
$ws.Range("M16").Value = 0.03004291845493562
$ws.Range("O16").Value = 0.03862660944206009
$ws.Range("S16").Value = 0.1888412017167382

# Row 17
$ws.Range("F17").Value = 0.03211009174311927
$ws.Range("H17").Value = 0.1513761467889908
$ws.Range("I17").Value = 0.1009174311926606
$ws.Range("J17").Value = 0.4541284403669725
$ws.Range("K17").Value = 0.0871559633027523
$ws.Range("M17").Value = 0.01605504587155963
$ws.Range("O17").Value = 0.06651376146788991
$ws.Range("S17").Value = 0.09174311926605505

# Row 18
$ws.Range("F18").Value = 0.03208556149732621
$ws.Range("H18").Value = 0.1925133689839572
$ws.Range("I18").Value = 0.1336898395721925
$ws.Range("J18").Value = 0.3957219251336899
$ws.Range("K18").Value = 0.0481283422459893
$ws.Range("M18").Value = 0.0160427807486631
$ws.Range("O18").Value = 0.0374331550802139
$ws.Range("S18").Value = 0.1443850267379679

# Row 19
$ws.Range("F19").Value = 0.02185792349726776
$ws.Range("H19").Value = 0.1943793911007026
$ws.Range("I19").Value = 0.07103825136612021
$ws.Range("J19").Value = 0.3965651834504293
$ws.Range("K19").Value = 0.1038251366120219
$ws.Range("M19").Value = 0.02029664324746292
$ws.Range("O19").Value = 0.06713505074160812
$ws.Range("S19").Value = 0.1249024199843872
